$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header strings (Volume number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"

# --- Update data table cells (rows 14-33) ---
$ws.Range("N14").Value = -81.081081081081
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = "***.*"
$ws.Range("L15").Value = -32.35294117647
$ws.Range("M15").Value = 4.545454545454
$ws.Range("N15").Value = -67.605633802816
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 119
$ws.Range("J16").Value = 122
$ws.Range("K16").Value = -2.459016393442
$ws.Range("L16").Value = 8.181818181818
$ws.Range("M16").Value = -59.661016949152
$ws.Range("N16").Value = -87.706611570247
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -37.931034482758
$ws.Range("I17").Value = 326
$ws.Range("J17").Value = 312
$ws.Range("K17").Value = 4.487179487179
$ws.Range("L17").Value = -2.97619047619
$ws.Range("M17").Value = 6.188925081433
$ws.Range("N17").Value = -50.229007633587
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 8
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 73
$ws.Range("J18").Value = 69
$ws.Range("K18").Value = 5.797101449275
$ws.Range("L18").Value = -34.234234234234
$ws.Range("M18").Value = -78.14371257485
$ws.Range("N18").Value = -92.364016736401
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = 16
$ws.Range("I19").Value = 264
$ws.Range("J19").Value = 346
$ws.Range("K19").Value = -23.699421965317
$ws.Range("L19").Value = -31.958762886597
$ws.Range("M19").Value = -53.356890459364
$ws.Range("N19").Value = -92.367736339982
$ws.Range("C20").Value = 4
$ws.Range("C16").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 5
$ws.Range("K14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 64.285714285714
$ws.Range("I20").Value = 201
$ws.Range("J20").Value = 186
$ws.Range("K20").Value = 8.064516129032
$ws.Range("L20").Value = -0.49504950495
$ws.Range("M20").Value = -14.468085106383
$ws.Range("N20").Value = -86.854153041203
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -7.692307692307
$ws.Range("F21").Value = 86
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = 1.176470588235
$ws.Range("I21").Value = 1013
$ws.Range("J21").Value = 1060
$ws.Range("K21").Value = -4.43396226415
$ws.Range("L21").Value = -14.8023549201
$ws.Range("M21").Value = -42.961711711711
$ws.Range("N21").Value = -86.801302931596
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("L23").Value = 16.666666666666
$ws.Range("M23").Value = -46.153846153846
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 13.333333333333
$ws.Range("F24").Value = 84
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = 2.439024390243
$ws.Range("I24").Value = 960
$ws.Range("J24").Value = 989
$ws.Range("K24").Value = -2.932254802831
$ws.Range("L24").Value = -16.376306620209
$ws.Range("M24").Value = 3.336921420882
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 26
$ws.Range("H25").Value = 8.333333333333
$ws.Range("I25").Value = 277
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = -7.666666666666
$ws.Range("L25").Value = -27.676240208877
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 19
$ws.Range("E26").Value = -15.78947368421
$ws.Range("G26").Value = 62
$ws.Range("H26").Value = 24.193548387096
$ws.Range("I26").Value = 714
$ws.Range("J26").Value = 585
$ws.Range("K26").Value = 22.051282051282
$ws.Range("L26").Value = 56.923076923076
$ws.Range("M26").Value = 0.847457627118
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = -15
$ws.Range("L27").Value = -29.166666666666
$ws.Range("C15").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = "***.*"
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 100
$ws.Range("C16").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1
$ws.Range("C15").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = "***.*"
$ws.Range("C16").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 16
$ws.Range("K29").Value = -11.111111111111
$ws.Range("L29").Value = -46.666666666666
$ws.Range("M29").Value = -68.627450980392
$ws.Range("N29").Value = -88.888888888888
$ws.Range("C16").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 1
$ws.Range("C15").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = "***.*"
$ws.Range("C16").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 14
$ws.Range("K30").Value = 7.692307692307
$ws.Range("L30").Value = -41.666666666666
$ws.Range("M30").Value = -66.666666666666
$ws.Range("N30").Value = -89.0625
$ws.Range("C15").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = "***.*"
$ws.Range("C16").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C33").Value = 1
$ws.Range("C16").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("F33").Value = 1
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 3
$ws.Range("K33").Value = -57.142857142857
$ws.Range("L33").Value = 0

$excel.CutCopyMode = 0

